$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "59÷9=6, 5"
$t.Cell(1, 2).Range.Text = "39÷4=9, 3"
$t.Cell(1, 3).Range.Text = "89÷7=12, 5"
$t.Cell(1, 4).Range.Text = "19÷5=3, 4"
$t.Cell(1, 5).Range.Text = "21÷6=3, 3"
$t.Cell(5, 1).Range.Text = "63÷4=15, 3"
$t.Cell(5, 2).Range.Text = "78÷5=15, 3"
$t.Cell(5, 3).Range.Text = "61÷9=6, 7"
$t.Cell(5, 4).Range.Text = "57÷5=11, 2"
$t.Cell(5, 5).Range.Text = "16÷2=8, 0"
$t.Cell(9, 1).Range.Text = "76÷2=38, 0"
$t.Cell(9, 2).Range.Text = "34÷2=17, 0"
$t.Cell(9, 3).Range.Text = "33÷2=16, 1"
$t.Cell(9, 4).Range.Text = "93÷6=15, 3"
$t.Cell(9, 5).Range.Text = "83÷8=10, 3"
$t.Cell(13, 1).Range.Text = "28÷5=5, 3"
$t.Cell(13, 2).Range.Text = "96÷8=12, 0"
$t.Cell(13, 3).Range.Text = "75÷4=18, 3"
$t.Cell(13, 4).Range.Text = "34÷7=4, 6"
$t.Cell(13, 5).Range.Text = "90÷2=45, 0"
$t.Cell(17, 1).Range.Text = "83÷3=27, 2"
$t.Cell(17, 2).Range.Text = "87÷9=9, 6"
$t.Cell(17, 3).Range.Text = "22÷2=11, 0"
$t.Cell(17, 4).Range.Text = "35÷5=7, 0"
$t.Cell(17, 5).Range.Text = "14÷3=4, 2"

Write-Output "Replacements applied"